$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 19).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 19)
    $val = $cell.Value2
    if ($val -ne $null -and $val.ToString().StartsWith("/img/")) {
        $cell.Value2 = $val.ToString().Substring(1)
    }
}
